{"js": "const replacements = [\n  [\"152\u00f77=21, 5\", \"228\u00f79=25, 3\"],\n  [\"505\u00f79=56, 1\", \"158\u00f77=22, 4\"],\n  [\"992\u00f77=141, 5\", \"831\u00f79=92, 3\"],\n  [\"123\u00f76=20, 3\", \"334\u00f77=47, 5\"],\n  [\"219\u00f74=54, 3\", \"739\u00f74=184, 3\"],\n  [\"516\u00f78=64, 4\", \"970\u00f76=161, 4\"],\n  [\"234\u00f74=58, 2\", \"182\u00f76=30, 2\"],\n  [\"238\u00f79=26, 4\", \"535\u00f75=107, 0\"],\n  [\"667\u00f74=166, 3\", \"315\u00f74=78, 3\"],\n  [\"757\u00f77=108, 1\", \"870\u00f78=108, 6\"],\n  [\"283\u00f73=94, 1\", \"150\u00f73=50, 0\"],\n  [\"103\u00f79=11, 4\", \"223\u00f75=44, 3\"],\n  [\"733\u00f72=366, 1\", \"619\u00f75=123, 4\"],\n  [\"778\u00f73=259, 1\", \"755\u00f74=188, 3\"],\n  [\"865\u00f75=173, 0\", \"423\u00f79=47, 0\"],\n  [\"943\u00f79=104, 7\", \"902\u00f79=100, 2\"],\n  [\"587\u00f73=195, 2\", \"496\u00f75=99, 1\"],\n  [\"290\u00f72=145, 0\", \"671\u00f79=74, 5\"],\n  [\"461\u00f79=51, 2\", \"290\u00f78=36, 2\"],\n  [\"692\u00f73=230, 2\", \"673\u00f78=84, 1\"],\n  [\"106\u00f75=21, 1\", \"271\u00f79=30, 1\"],\n  [\"397\u00f76=66, 1\", \"564\u00f72=282, 0\"],\n  [\"839\u00f73=279, 2\", \"918\u00f75=183, 3\"],\n  [\"643\u00f77=91, 6\", \"324\u00f79=36, 0\"],\n  [\"778\u00f75=155, 3\", \"895\u00f73=298, 1\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();", "ps1": "$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"152\u00f77=21, 5\", \"228\u00f79=25, 3\"),\n  @(\"505\u00f79=56, 1\", \"158\u00f77=22, 4\"),\n  @(\"992\u00f77=141, 5\", \"831\u00f79=92, 3\"),\n  @(\"123\u00f76=20, 3\", \"334\u00f77=47, 5\"),\n  @(\"219\u00f74=54, 3\", \"739\u00f74=184, 3\"),\n  @(\"516\u00f78=64, 4\", \"970\u00f76=161, 4\"),\n  @(\"234\u00f74=58, 2\", \"182\u00f76=30, 2\"),\n  @(\"238\u00f79=26, 4\", \"535\u00f75=107, 0\"),\n  @(\"667\u00f74=166, 3\", \"315\u00f74=78, 3\"),\n  @(\"757\u00f77=108, 1\", \"870\u00f78=108, 6\"),\n  @(\"283\u00f73=94, 1\", \"150\u00f73=50, 0\"),\n  @(\"103\u00f79=11, 4\", \"223\u00f75=44, 3\"),\n  @(\"733\u00f72=366, 1\", \"619\u00f75=123, 4\"),\n  @(\"778\u00f73=259, 1\", \"755\u00f74=188, 3\"),\n  @(\"865\u00f75=173, 0\", \"423\u00f79=47, 0\"),\n  @(\"943\u00f79=104, 7\", \"902\u00f79=100, 2\"),\n  @(\"587\u00f73=195, 2\", \"496\u00f75=99, 1\"),\n  @(\"290\u00f72=145, 0\", \"671\u00f79=74, 5\"),\n  @(\"461\u00f79=51, 2\", \"290\u00f78=36, 2\"),\n  @(\"692\u00f73=230, 2\", \"673\u00f78=84, 1\"),\n  @(\"106\u00f75=21, 1\", \"271\u00f79=30, 1\"),\n  @(\"397\u00f76=66, 1\", \"564\u00f72=282, 0\"),\n  @(\"839\u00f73=279, 2\", \"918\u00f75=183, 3\"),\n  @(\"643\u00f77=91, 6\", \"324\u00f79=36, 0\"),\n  @(\"778\u00f75=155, 3\", \"895\u00f73=298, 1\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute([ref]$find.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]$wdFindContinue, [ref]$false, [ref]$find.Replacement.Text, [ref]$wdReplaceAll)\n}"}
